# Auto-generated edit script applying the cryptos.xlsx diff
# Updates Price (D) and Volume(1h) (E) columns for each row, and
# swaps the BitcoinCash/Uniswap (rows 20-21) and Binance-PegBSC-USD/PEPE (rows 29-30) entries.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.202.44"
$ws.Range("E2").Value = "  -2.09%  "
$ws.Range("D3").Value = "2.664.06"
$ws.Range("E3").Value = "  -1.46%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "597.53"
$ws.Range("E5").Value = "  -0.29%  "
$ws.Range("D6").Value = "165.56"
$ws.Range("E6").Value = "  +3.19%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").Value = "0.547"
$ws.Range("E8").Value = "  +0.36%  "
$ws.Range("D9").Value = "2.663.54"
$ws.Range("E9").Value = "  -1.45%  "
$ws.Range("E10").Value = "  +0.56%  "
$ws.Range("E11").Value = "  +1.22%  "
$ws.Range("E12").Value = "  -1.00%  "
$ws.Range("D13").Value = "5.22"
$ws.Range("E13").Value = "  -1.75%  "
$ws.Range("D14").Value = "27.72"
$ws.Range("E14").Value = "  -2.17%  "
$ws.Range("D15").Value = "3.149.65"
$ws.Range("E15").Value = "  -1.46%  "
$ws.Range("D16").Value = "0.0000184"
$ws.Range("E16").Value = "  -2.86%  "
$ws.Range("D17").Value = "67.101.70"
$ws.Range("E17").Value = "  -2.13%  "
$ws.Range("D18").Value = "2.665.26"
$ws.Range("E18").Value = "  -1.34%  "
$ws.Range("D19").Value = "11.69"
$ws.Range("E19").Value = "  -1.60%  "
$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").Value = "7.61"
$ws.Range("E20").Value = "  -0.27%  "
$ws.Range("B21").Value = "BitcoinCash"
$ws.Range("C21").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D21").Value = "362.58"
$ws.Range("E21").Value = "  -1.00%  "
$ws.Range("E22").Value = "  -3.74%  "
$ws.Range("E23").Value = "  -2.25%  "
$ws.Range("E24").Value = "  -4.88%  "
$ws.Range("E25").Value = "  +0.02%  "
$ws.Range("D26").Value = "70.78"
$ws.Range("E26").Value = "  -5.02%  "
$ws.Range("D27").Value = "10.04"
$ws.Range("E27").Value = "  +0.62%  "
$ws.Range("D28").Value = "2.823.96"
$ws.Range("E28").Value = "  -0.12%  "
$ws.Range("B29").Value = "PEPE"
$ws.Range("C29").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D29").Value = "0.0000102"
$ws.Range("E29").Value = "  -3.69%  "
$ws.Range("B30").Value = "Binance-PegBSC-USD"
$ws.Range("C30").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D30").Value = "1.00"
$ws.Range("E30").Value = "  +0.08%  "
$ws.Range("D31").Value = "555.04"
$ws.Range("E31").Value = "  -4.12%  "
$ws.Range("D32").Value = "7.99"
$ws.Range("E32").Value = "  -3.26%  "
$ws.Range("E33").Value = "  -4.82%  "
$ws.Range("E34").Value = "  -1.78%  "
$ws.Range("E35").Value = "  -2.45%  "
$ws.Range("E36").Value = "  -0.01%  "
$ws.Range("E37").Value = "  -5.82%  "
$ws.Range("E38").Value = "  -2.09%  "
$ws.Range("D39").Value = "156.13"
$ws.Range("E39").Value = "  -2.82%  "
$ws.Range("E40").Value = "  -2.09%  "
$ws.Range("D41").Value = "5.29"
$ws.Range("E41").Value = "  -2.32%  "
$ws.Range("D42").Value = "1.83"
$ws.Range("E42").Value = "  -4.74%  "
$ws.Range("E43").Value = "  +0.39%  "
$ws.Range("E44").Value = "  +0.00%  "
$ws.Range("E45").Value = "  -7.30%  "
$ws.Range("D46").Value = "40.32"
$ws.Range("E46").Value = "  -0.81%  "
$ws.Range("E47").Value = "  -6.98%  "
$ws.Range("E48").Value = "  -2.22%  "
$ws.Range("D49").Value = "152.66"
$ws.Range("E49").Value = "  -3.41%  "
$ws.Range("E50").Value = "  -3.50%  "
$ws.Range("E51").Value = "  -3.60%  "
